$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Add new row 7 with the new mail log entry
$ws.Range("A7").Value = "Nieuwe bestelling"
$ws.Range("B7").Value = "planning@testbedrijf123.nl"
$ws.Range("D7").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("F7").Value = "2025-08-19 19:41:24"
$ws.Range("G7").Value = "Nee"
$ws.Range("H7").Value = "Ja"
$ws.Range("I7").Value = "Nee"
$ws.Range("J7").Value = "Nee"

# Extend conditional formatting ranges to include the new row (D/G/H/I/J 2:6 -> 2:7)
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $ws.Range($col + "2:" + $col + "6")
    $newRange = $ws.Range($col + "2:" + $col + "7")
    foreach ($fc in $oldRange.FormatConditions) {
        $fc.ModifyAppliesToRange($newRange)
    }
}

# Update Dashboard count
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 6
